$wb = $excel.ActiveWorkbook

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1023.375
$ws.Range("J129").Value = 1115
$ws.Range("L129").Value = 3345
$ws.Range("N129").Value = -13345

# ALC row 131
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 1887.4032
$ws.Range("I131").Value = 392.18182
$ws.Range("J131").Value = 2209.9019
$ws.Range("K131").Value = 1176.54546
$ws.Range("L131").Value = 6629.7057
$ws.Range("M131").Value = 3863.45454
$ws.Range("N131").Value = -16709.7057

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1914.5834
$ws.Range("I141").Value = 1999
$ws.Range("K141").Value = 5997
$ws.Range("M141").Value = -817

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5299.2607
$ws.Range("I2").Value = 1003.85
$ws.Range("J2").Value = 33935.332
$ws.Range("K2").Value = 1003.85
$ws.Range("L2").Value = 33935.332
$ws.Range("M2").Value = -890.85
$ws.Range("N2").Value = -34161.332

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21899.549
$ws.Range("I32").Value = 14886.185
$ws.Range("J32").Value = 38180.57
$ws.Range("K32").Value = 14886.185
$ws.Range("L32").Value = 38180.57
$ws.Range("M32").Value = -14599.185
$ws.Range("N32").Value = -38754.57

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 19760.926
$ws.Range("I74").Value = 22801.436
$ws.Range("J74").Value = 2278
$ws.Range("K74").Value = 22801.436
$ws.Range("L74").Value = 2278
$ws.Range("M74").Value = -21927.436
$ws.Range("N74").Value = -4026

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 19760.926
$ws.Range("I77").Value = 22801.436
$ws.Range("J77").Value = 2278
$ws.Range("K77").Value = 114007.18
$ws.Range("L77").Value = 11390
$ws.Range("M77").Value = -109639.18
$ws.Range("N77").Value = -20126

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 5299.2607
$ws.Range("I116").Value = 1003.85
$ws.Range("J116").Value = 33935.332
$ws.Range("K116").Value = 1003.85
$ws.Range("L116").Value = 33935.332
$ws.Range("M116").Value = 1290.15
$ws.Range("N116").Value = -38523.332

# ARM row 124
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H124").Value = 11066.667
$ws.Range("J124").Value = 11066.667
$ws.Range("L124").Value = 11066.667
$ws.Range("N124").Value = -20886.667

# ARM row 125
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H125").Value = 33750
$ws.Range("J125").Value = 33750
$ws.Range("L125").Value = 33750
$ws.Range("N125").Value = -43590

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 221600.14
$ws.Range("I132").Value = 32387.05
$ws.Range("J132").Value = 836542.7
$ws.Range("K132").Value = 97161.14999999999
$ws.Range("L132").Value = 2509628.1
$ws.Range("M132").Value = -94631.14999999999
$ws.Range("N132").Value = -2514688.1

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5299.2607
$ws.Range("I3").Value = 1003.85
$ws.Range("J3").Value = 33935.332
$ws.Range("K3").Value = 1003.85
$ws.Range("L3").Value = 33935.332
$ws.Range("M3").Value = -889.85
$ws.Range("N3").Value = -34163.332

# BSM row 74
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 33390
$ws.Range("J74").Value = 33390
$ws.Range("L74").Value = 33390
$ws.Range("N74").Value = -35262

# BSM row 77
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H77").Value = 33390
$ws.Range("J77").Value = 33390
$ws.Range("L77").Value = 100170
$ws.Range("N77").Value = -109530

# CRP row 17
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 34006.5
$ws.Range("I17").Value = 8004
$ws.Range("K17").Value = 8004
$ws.Range("M17").Value = -7830

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 30106.645
$ws.Range("I31").Value = 35229.9
$ws.Range("J31").Value = 19860.133
$ws.Range("K31").Value = 35229.9
$ws.Range("L31").Value = 19860.133
$ws.Range("M31").Value = -34934.9
$ws.Range("N31").Value = -20450.133

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 30106.645
$ws.Range("I34").Value = 35229.9
$ws.Range("J34").Value = 19860.133
$ws.Range("K34").Value = 35229.9
$ws.Range("L34").Value = 19860.133
$ws.Range("M34").Value = -35027.9
$ws.Range("N34").Value = -20264.133

# CUL row 75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 3026.875
$ws.Range("I75").Value = 200
$ws.Range("J75").Value = 3283.8635
$ws.Range("K75").Value = 600
$ws.Range("L75").Value = 9851.5905
$ws.Range("M75").Value = 398
$ws.Range("N75").Value = -11847.5905

# CUL row 78
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 3026.875
$ws.Range("I78").Value = 200
$ws.Range("J78").Value = 3283.8635
$ws.Range("K78").Value = 1800
$ws.Range("L78").Value = 29554.7715
$ws.Range("M78").Value = 3192
$ws.Range("N78").Value = -39538.7715

# CUL row 114
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 404
$ws.Range("I114").Value = 290.18182
$ws.Range("J114").Value = 529.2
$ws.Range("K114").Value = 870.54546
$ws.Range("L114").Value = 1587.6
$ws.Range("M114").Value = 2383.45454
$ws.Range("N114").Value = -8095.6

# CUL row 117
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 1828.5
$ws.Range("I117").Value = 500
$ws.Range("J117").Value = 2018.2858
$ws.Range("K117").Value = 1500
$ws.Range("L117").Value = 6054.857400000001
$ws.Range("M117").Value = 1942
$ws.Range("N117").Value = -12938.8574

# CUL row 121
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 62562910
$ws.Range("I121").Value = 400
$ws.Range("J121").Value = 100100424
$ws.Range("K121").Value = 1200
$ws.Range("L121").Value = 300301272
$ws.Range("M121").Value = 110
$ws.Range("N121").Value = -300303892

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 176272.61
$ws.Range("I131").Value = 445.44446
$ws.Range("J131").Value = 209240.2
$ws.Range("K131").Value = 1336.33338
$ws.Range("L131").Value = 627720.6000000001
$ws.Range("M131").Value = 3703.66662
$ws.Range("N131").Value = -637800.6000000001

# CUL row 137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 28818718
$ws.Range("I137").Value = 111112376
$ws.Range("J137").Value = 5306246
$ws.Range("K137").Value = 333337128
$ws.Range("L137").Value = 15918738
$ws.Range("M137").Value = -333332028
$ws.Range("N137").Value = -15928938

# GSM row 96
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1844.2858
$ws.Range("I126").Value = 1500
$ws.Range("K126").Value = 4500
$ws.Range("M126").Value = -2030

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2573.075
$ws.Range("I122").Value = 2751.913
$ws.Range("J122").Value = 2331.1177
$ws.Range("K122").Value = 8255.739
$ws.Range("L122").Value = 6993.353099999999
$ws.Range("M122").Value = -5805.739
$ws.Range("N122").Value = -11893.3531

# WVR row 141
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 52678.184
$ws.Range("J141").Value = 52678.184
$ws.Range("L141").Value = 52678.184
$ws.Range("N141").Value = -63038.184
